# Aanpassing verslag + dailyscrum
$wb = $excel.ActiveWorkbook

# --- 02-03-2016: fill in the remaining "problems" column (G) for rows 3 and 4 ---
$ws3 = $wb.Worksheets.Item("02-03-2016")
$ws3.Range("G3").Value = "Use case maken over het filteringsysteem en mutual friends, een paginabeschrijving over het filteringsysteem en het de Facebook API bespreken"
$ws3.Range("G4").Value = "Afmaken Facebook API + bundelen paginabeschrijvingen"

# --- update active-cell selections left behind by the editing session ---
$ws2 = $wb.Worksheets.Item("26-02-2016")
[void]$ws2.Range("G4").Select()

[void]$ws3.Range("F4").Select()
